$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 954.9545000000001
$ws.Range("I39").Value = 954.9545000000001
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2864.8635
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -2568.8635

$ws.Range("H42").Value = 329.5
$ws.Range("I42").Value = 88.666664
$ws.Range("J42").Value = 474
$ws.Range("K42").Value = 265.999992
$ws.Range("L42").Value = 1422
$ws.Range("M42").Value = -35.99999200000002
$ws.Range("N42").Value = -1882

$ws.Range("I92").Value = 1201
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1201
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = 47

$ws.Range("H112").Value = 2525.8696
$ws.Range("J112").Value = 2570.158
$ws.Range("L112").Value = 7710.474
$ws.Range("N112").Value = -9926.474

$ws.Range("H125").Value = 41452.5
$ws.Range("I125").Value = 65395.6
$ws.Range("K125").Value = 588560.4
$ws.Range("M125").Value = -586100.4

$ws.Range("H127").Value = 92169.89999999999
$ws.Range("I127").Value = 92169.89999999999
$ws.Range("K127").Value = 276509.7
$ws.Range("M127").Value = -271549.7

$ws.Range("H132").Value = 2700.6042
$ws.Range("I132").Value = 1069.2703
$ws.Range("J132").Value = 8187.8184
$ws.Range("K132").Value = 3207.810899999999
$ws.Range("L132").Value = 24563.4552
$ws.Range("M132").Value = -677.8108999999995
$ws.Range("N132").Value = -29623.4552

$ws.Range("H137").Value = 3701.561
$ws.Range("I137").Value = 1375.4667
$ws.Range("J137").Value = 10045.454
$ws.Range("K137").Value = 4126.4001
$ws.Range("L137").Value = 30136.362
$ws.Range("M137").Value = -1576.4001
$ws.Range("N137").Value = -35236.362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = 0

$ws.Range("H32").Value = 34026.125
$ws.Range("J32").Value = 158221
$ws.Range("L32").Value = 158221
$ws.Range("N32").Value = -158795

$ws.Range("H132").Value = 12001.737
$ws.Range("I132").Value = 12913.706
$ws.Range("K132").Value = 38741.118
$ws.Range("M132").Value = -36211.118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 369.8
$ws.Range("I22").Value = 369.8
$ws.Range("K22").Value = 369.8
$ws.Range("M22").Value = -196.8

$ws.Range("H24").Value = 3758
$ws.Range("I24").Value = 3758
$ws.Range("K24").Value = 3758
$ws.Range("M24").Value = -3523

$ws.Range("H29").Value = 1758
$ws.Range("I29").Value = 1758
$ws.Range("K29").Value = 1758
$ws.Range("M29").Value = -1469

$ws.Range("H35").Value = 29673.334
$ws.Range("I35").Value = 100
$ws.Range("K35").Value = 100
$ws.Range("M35").Value = 210

$ws.Range("H86").Value = 4082.4443
$ws.Range("I86").Value = 3804.6667
$ws.Range("K86").Value = 3804.6667
$ws.Range("M86").Value = -2681.6667

$ws.Range("H89").Value = 4082.4443
$ws.Range("I89").Value = 3804.6667
$ws.Range("K89").Value = 19023.3335
$ws.Range("M89").Value = -13407.3335

$ws.Range("H134").Value = 1292.2
$ws.Range("I134").Value = 1295.091
$ws.Range("J134").Value = 1271
$ws.Range("K134").Value = 3885.273
$ws.Range("L134").Value = 3813
$ws.Range("M134").Value = -1350.273
$ws.Range("N134").Value = -8883

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1200.8182
$ws.Range("I16").Value = 926.375
$ws.Range("J16").Value = 1932.6666
$ws.Range("K16").Value = 926.375
$ws.Range("L16").Value = 1932.6666
$ws.Range("M16").Value = -639.375
$ws.Range("N16").Value = -2506.6666

$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 192.3077
$ws.Range("K22").Value = 192.3077
$ws.Range("M22").Value = 157.6923

$ws.Range("H62").Value = 3899.8572
$ws.Range("I62").Value = 3899.8572
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3899.8572
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -3275.8572

$ws.Range("H65").Value = 3899.8572
$ws.Range("I65").Value = 3899.8572
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 19499.286
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -16379.286

$ws.Range("H99").Value = 1806.6923
$ws.Range("J99").Value = 2120
$ws.Range("L99").Value = 2120
$ws.Range("N99").Value = -5116

$ws.Range("H113").Value = 1200.8182
$ws.Range("I113").Value = 926.375
$ws.Range("J113").Value = 1932.6666
$ws.Range("K113").Value = 926.375
$ws.Range("L113").Value = 1932.6666
$ws.Range("M113").Value = 1243.625
$ws.Range("N113").Value = -6272.6666

$ws.Range("H126").Value = 1806.6923
$ws.Range("J126").Value = 2120
$ws.Range("L126").Value = 6360
$ws.Range("N126").Value = -11300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 667.375
$ws.Range("I5").Value = 655
$ws.Range("J5").Value = 704.5
$ws.Range("K5").Value = 1965
$ws.Range("L5").Value = 2113.5
$ws.Range("M5").Value = -1853
$ws.Range("N5").Value = -2337.5

$ws.Range("H135").Value = 667.375
$ws.Range("I135").Value = 655
$ws.Range("J135").Value = 704.5
$ws.Range("K135").Value = 5895
$ws.Range("L135").Value = 6340.5
$ws.Range("M135").Value = -3360
$ws.Range("N135").Value = -11410.5

$ws.Range("H137").Value = 4549191
$ws.Range("I137").Value = 11112626
$ws.Range("J137").Value = 5274.5386
$ws.Range("K137").Value = 33337878
$ws.Range("L137").Value = 15823.6158
$ws.Range("M137").Value = -33332778
$ws.Range("N137").Value = -26023.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 83292.625
$ws.Range("I97").Value = 94334.42999999999
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 94334.42999999999
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -93838.42999999999
$ws.Range("N97").Value = -6992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6048.1665
$ws.Range("I40").Value = 4131.3335
$ws.Range("J40").Value = 7965
$ws.Range("K40").Value = 4131.3335
$ws.Range("L40").Value = 7965
$ws.Range("M40").Value = -3995.3335
$ws.Range("N40").Value = -8237

$ws.Range("H68").Value = 2312.8333
$ws.Range("I68").Value = 2225.7144
$ws.Range("J68").Value = 2434.8
$ws.Range("K68").Value = 2225.7144
$ws.Range("L68").Value = 2434.8
$ws.Range("M68").Value = -1476.7144
$ws.Range("N68").Value = -3932.8

$ws.Range("H71").Value = 2312.8333
$ws.Range("I71").Value = 2225.7144
$ws.Range("J71").Value = 2434.8
$ws.Range("K71").Value = 11128.572
$ws.Range("L71").Value = 12174
$ws.Range("M71").Value = -7384.572
$ws.Range("N71").Value = -19662

$ws.Range("H82").Value = 1359.2307
$ws.Range("I82").Value = 1725
$ws.Range("J82").Value = 1196.6666
$ws.Range("K82").Value = 1725
$ws.Range("L82").Value = 1196.6666
$ws.Range("M82").Value = -1364
$ws.Range("N82").Value = -1918.6666

$ws.Range("H85").Value = 1359.2307
$ws.Range("I85").Value = 1725
$ws.Range("J85").Value = 1196.6666
$ws.Range("K85").Value = 1725
$ws.Range("L85").Value = 1196.6666
$ws.Range("M85").Value = -477
$ws.Range("N85").Value = -3692.6666

$ws.Range("H93").Value = 32812.91
$ws.Range("I93").Value = 2623.2222
$ws.Range("K93").Value = 2623.2222
$ws.Range("M93").Value = -1375.2222

$ws.Range("H132").Value = 2688.1516
$ws.Range("I132").Value = 2190.4666
$ws.Range("K132").Value = 6571.399800000001
$ws.Range("M132").Value = -4041.399800000001

$ws.Range("H136").Value = 2851.9048
$ws.Range("I136").Value = 2594.75
$ws.Range("J136").Value = 7995
$ws.Range("K136").Value = 7784.25
$ws.Range("L136").Value = 23985
$ws.Range("M136").Value = -5234.25
$ws.Range("N136").Value = -29085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51000
$ws.Range("J46").Value = 51000
$ws.Range("L46").Value = 51000
$ws.Range("N46").Value = -51462

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = ""
$ws.Range("N51").Value = 0

$ws.Range("H52").Value = 18856.285
$ws.Range("J52").Value = 20332.5
$ws.Range("L52").Value = 20332.5
$ws.Range("N52").Value = -20784.5

$ws.Range("H100").Value = 3552.3845
$ws.Range("I100").Value = 3848.2
$ws.Range("K100").Value = 7696.4
$ws.Range("M100").Value = -7155.4

$ws.Range("H134").Value = 51000
$ws.Range("J134").Value = 51000
$ws.Range("L134").Value = 153000
$ws.Range("N134").Value = -158070

$ws.Range("H136").Value = 1869.279
$ws.Range("I136").Value = 1869.279
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5607.837
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -3057.837
